$wb = $excel.ActiveWorkbook
$wbWs = $wb.Worksheets.Item("Workblocks")
$wbWs.Delete()

$introWs = $wb.Worksheets.Item("Introduction")
$introWs.Rows("8:9").Delete()
